# Apply updated loading-time values after modifying penalties/fixing
# assignments and new population generation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = 5

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 10

# Row 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 1

# Row 5
$ws.Range("B5").Value = 2
$ws.Range("D5").Value = 7

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5

# Row 7
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 10
